$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date value of 45203 for every data row
# (rows 2-353). Update it to 45205 for all of them.
$ws.Range("C2:C353").Value = 45205
